# Updated cryptos list on Sun Aug 27 14:10:50 UTC 2023 with GitHub Actions
#
# Applies the per-cell Price/Volume(1h) refresh (and the Frax <-> BabyDogeCoin
# row swap) from the upstream scrape run. Cells are addressed by A1 ref on
# Sheet1 (the workbook's ActiveSheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell this commit touches, as (A1 ref, new text) pairs. Values that
# look like a bare number (e.g. "0.8660", "1.011") are written with a
# leading apostrophe so Excel's COM layer stores them as literal text
# instead of silently re-parsing them into a Double (which would drop
# meaningful trailing zeros such as 0.8660 -> 0.866). Everything else
# (two-dot "thousands" prices like "26.367.80", the padded "  +1.11%  "
# strings, coin names, and coinranking.com URLs) is never numeric-looking,
# so it round-trips through .Value untouched.
$changes = @(
    @{ Cell = "D2"; Value = "26.367.80" },
    @{ Cell = "E2"; Value = "  +1.11%  " },
    @{ Cell = "D3"; Value = "1.669.45" },
    @{ Cell = "E3"; Value = "  +1.06%  " },
    @{ Cell = "D5"; Value = "219.92" },
    @{ Cell = "E5"; Value = "  +1.14%  " },
    @{ Cell = "D6"; Value = "0.5351" },
    @{ Cell = "E6"; Value = "  +1.29%  " },
    @{ Cell = "E7"; Value = "  +0.93%  " },
    @{ Cell = "D8"; Value = "0.2657" },
    @{ Cell = "E8"; Value = "  +2.25%  " },
    @{ Cell = "D9"; Value = "0.06396" },
    @{ Cell = "E9"; Value = "  +1.31%  " },
    @{ Cell = "D10"; Value = "20.95" },
    @{ Cell = "E10"; Value = "  +3.01%  " },
    @{ Cell = "D11"; Value = "0.07868" },
    @{ Cell = "E11"; Value = "  +1.01%  " },
    @{ Cell = "D12"; Value = "4.567" },
    @{ Cell = "E12"; Value = "  +1.09%  " },
    @{ Cell = "D13"; Value = "1.673.93" },
    @{ Cell = "E13"; Value = "  +1.17%  " },
    @{ Cell = "D14"; Value = "1.897.82" },
    @{ Cell = "E14"; Value = "  +1.01%  " },
    @{ Cell = "D15"; Value = "0.5546" },
    @{ Cell = "E15"; Value = "  +1.26%  " },
    @{ Cell = "D16"; Value = "0.0₅8169" },
    @{ Cell = "E16"; Value = "  -0.37%  " },
    @{ Cell = "D17"; Value = "66.17" },
    @{ Cell = "E17"; Value = "  +1.34%  " },
    @{ Cell = "D18"; Value = "26.399.75" },
    @{ Cell = "E18"; Value = "  +1.23%  " },
    @{ Cell = "E19"; Value = "  +0.94%  " },
    @{ Cell = "E20"; Value = "  +1.99%  " },
    @{ Cell = "D21"; Value = "196.10" },
    @{ Cell = "E21"; Value = "  +2.78%  " },
    @{ Cell = "D22"; Value = "10.28" },
    @{ Cell = "E22"; Value = "  +2.24%  " },
    @{ Cell = "D23"; Value = "6.046" },
    @{ Cell = "E23"; Value = "  +0.51%  " },
    @{ Cell = "E24"; Value = "  +0.95%  " },
    @{ Cell = "D25"; Value = "145.99" },
    @{ Cell = "E25"; Value = "  +1.38%  " },
    @{ Cell = "E26"; Value = "  -0.34%  " },
    @{ Cell = "D27"; Value = "7.249" },
    @{ Cell = "E27"; Value = "  +0.55%  " },
    @{ Cell = "D28"; Value = "16.17" },
    @{ Cell = "E28"; Value = "  +1.14%  " },
    @{ Cell = "D29"; Value = "1.502" },
    @{ Cell = "E29"; Value = "  +3.41%  " },
    @{ Cell = "D30"; Value = "0.05865" },
    @{ Cell = "E30"; Value = "  +1.40%  " },
    @{ Cell = "D31"; Value = "1.287" },
    @{ Cell = "E31"; Value = "  +1.28%  " },
    @{ Cell = "D32"; Value = "3.578" },
    @{ Cell = "E32"; Value = "  +0.93%  " },
    @{ Cell = "D33"; Value = "3.299" },
    @{ Cell = "E33"; Value = "  +1.24%  " },
    @{ Cell = "D34"; Value = "1.616" },
    @{ Cell = "E34"; Value = "  +1.31%  " },
    @{ Cell = "D35"; Value = "0.9709" },
    @{ Cell = "E35"; Value = "  +2.83%  " },
    @{ Cell = "D36"; Value = "2.838" },
    @{ Cell = "E36"; Value = "  +1.64%  " },
    @{ Cell = "D37"; Value = "2.432" },
    @{ Cell = "E37"; Value = "  +0.82%  " },
    @{ Cell = "D38"; Value = "0.5816" },
    @{ Cell = "E38"; Value = "  +1.27%  " },
    @{ Cell = "D39"; Value = "0.01609" },
    @{ Cell = "E39"; Value = "  -0.27%  " },
    @{ Cell = "D40"; Value = "1.075.98" },
    @{ Cell = "E40"; Value = "  +4.37%  " },
    @{ Cell = "D41"; Value = "0.8660" },
    @{ Cell = "E41"; Value = "  +1.43%  " },
    @{ Cell = "D42"; Value = "5.870" },
    @{ Cell = "E42"; Value = "  +2.87%  " },
    @{ Cell = "D43"; Value = "1.011" },
    @{ Cell = "E43"; Value = "  +0.99%  " },
    @{ Cell = "D44"; Value = "104.30" },
    @{ Cell = "E44"; Value = "  +0.15%  " },
    @{ Cell = "D45"; Value = "1.807.00" },
    @{ Cell = "E45"; Value = "  +0.74%  " },
    @{ Cell = "D46"; Value = "58.13" },
    @{ Cell = "E46"; Value = "  +2.21%  " },
    @{ Cell = "B47"; Value = "BabyDogeCoin" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" },
    @{ Cell = "D47"; Value = "0.0₈107" },
    @{ Cell = "E47"; Value = "  -4.13%  " },
    @{ Cell = "B48"; Value = "Frax" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax" },
    @{ Cell = "D48"; Value = "1.016" },
    @{ Cell = "E48"; Value = "  +1.51%  " },
    @{ Cell = "D49"; Value = "0.4394" },
    @{ Cell = "E49"; Value = "  +1.55%  " },
    @{ Cell = "D50"; Value = "8.070" },
    @{ Cell = "E50"; Value = "  +2.65%  " },
    @{ Cell = "E51"; Value = "  +0.51%  " }
)

$numberPattern = '^[+-]?\d+(\.\d+)?$'

foreach ($change in $changes) {
    $value = $change.Value
    if ($value -match $numberPattern) {
        $value = "'" + $value
    }
    $ws.Range($change.Cell).Value = $value
}
